$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells, re-using the existing header formatting
# (bold, centered, bordered) by copying the format from the adjacent
# header cell H1, then setting the text.
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new I0 / IF columns, rows 2-32
$data = @(
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(8, 9),
    @(8, 8),
    @(8, 8),
    @(10, 11),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(6, 6),
    @(7, 7),
    @(6, 7),
    @(6, 6),
    @(8, 8),
    @(8, 8),
    @(6, 7),
    @(8, 8),
    @(7, 8),
    @(10, 10),
    @(9, 9),
    @(8, 9),
    @(9, 9),
    @(9, 9),
    @(6, 7),
    @(8, 8),
    @(9, 9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
